$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: BU_B line 4
$ws.Cells.Item(5, 1).Value = 40500
$ws.Cells.Item(5, 2).Value = "BRL"
$ws.Cells.Item(5, 3).Value = "ACC-BU_B-0004"
$ws.Cells.Item(5, 4).Value = "Sample closure line 4 for BU_B"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "2025-02-15"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 6).Value = "BU_B"

# Row 6: BU_B line 5
$ws.Cells.Item(6, 1).Value = 50500
$ws.Cells.Item(6, 2).Value = "BRL"
$ws.Cells.Item(6, 3).Value = "ACC-BU_B-0005"
$ws.Cells.Item(6, 4).Value = "Sample closure line 5 for BU_B"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "2025-02-15"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 6).Value = "BU_B"
